# Updated cryptos list on Sat May 18 16:56:39 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.908.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.19%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.112.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.52%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.00%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.109.61"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.46%  "

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.75%  "

# Row 10 - Toncoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.04%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -1.35%  "

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.484"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.45%  "

# Row 13 - ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000245"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.74%  "

# Row 14 - Avalanche
$ws.Range("E14").Value = "  +0.87%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  -1.33%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.629.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.36%  "

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.867.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.29%  "

# Row 19 - WrappedEther
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.113.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.42%  "

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.53%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "476.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.46%  "

# Row 22 - Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.714"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.28%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +4.97%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("E24").Value = "  +4.50%  "

# Row 25 - Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.22%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  -2.01%  "

# Row 27 - RenderToken
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.19%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  -0.05%  "

# Row 29 - NEARProtocol
$ws.Range("E29").Value = "  -2.22%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  -2.04%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.27%  "

# Row 32 - EthereumClassic
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.60%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  +0.25%  "

# Row 34 - PEPE
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0939"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.26%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  +0.01%  "

# Row 36 - Filecoin
$ws.Range("E36").Value = "  -1.26%  "

# Row 37 - Mantle
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.978"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.37%  "

# Row 38 - Arweave
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.48%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  -0.83%  "

# Row 40 - OKB
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.69%  "

# Row 41 - TheGraph
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.311"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.61%  "

# Row 42 - Kaspa
$ws.Range("E42").Value = "  -1.25%  "

# Row 43 - Cosmos
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.57%  "

# Row 44 - Maker
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.800.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.32%  "

# Row 45
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.17%  "

# Row 46
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0356"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.02%  "

# Row 47 - Bittensor
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "380.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.48%  "

# Row 48 - Monero
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.59%  "

# Row 49 - USDe
$ws.Range("E49").Value = "  +0.06%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.05%  "

# Row 51 - ThetaToken
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.40%  "
